$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New temperature (column B) values for rows 2..166 (data index 0..164)
$bValues = @(218.9504919815064,162.9428252887726,210.9375297355652,156.9407932186127,147.5167701911927,169.7609234714508,161.0787859535218,135.7958772563935,191.7397353839875,165.1139590740204,156.382120552063,157.4156947898865,171.1493061828613,171.6427793979645,179.573070011139,218.0198698234558,222.7439555549622,183.9886836528779,146.9267032337189,212.6245820999146,198.2250938224793,193.6845123577118,154.7787862968445,161.688524017334,186.6197106838226,152.8565416812897,137.8981158065797,189.1655359935761,196.4158951473237,191.5680867671967,118.1949131107331,220.319547996521,205.1891953277589,146.0390660858155,222.883123664856,213.4426351928711,147.9150273036958,206.3845362663269,206.3006327056885,153.733621339798,180.4826434230805,168.7729305648804,199.3803374958039,205.6767208480835,216.4421808624268,219.7417340850831,192.234996881485,194.1080041217804,219.017254524231,186.5459539699555,169.7810714530945,237.4104410171509,226.1460899543763,144.3581203746796,118.5343688488007,199.3860997390748,149.0802191257477,141.7160589790345,144.6504184436799,143.0428188037873,209.1268342018128,185.3883789062501,183.3768261432648,156.2031466007233,234.1883800315857,186.519010515213,167.0206257629395,128.9151362037659,176.5374142837525,173.4187292671204,170.3716815185547,107.6412784576416,208.8461136054993,155.1084263515473,184.2105101490021,150.6782812595368,173.3951504325867,173.9964636993409,165.6525367450714,158.5893908882142,180.2994173336029,197.1385996723176,148.9178298377991,159.6959667873383,187.5396296405793,207.328855342865,173.6854482746124,167.8326516819001,182.7627034759522,169.1217780971527,202.316193819046,171.9082664775849,213.530512714386,198.226325750351,181.8717811870575,222.7012487220765,197.3377348518372,215.1444706916809,114.2310328483582,211.90672580719,166.3111545276643,157.6827582073212,208.5590346717835,202.541795578003,172.0227430438996,157.5648640346527,195.1310208702088,218.8415789604187,190.9199602413178,167.9628916263581,160.3895289325715,161.2865578746796,172.0795044517517,223.3589127540589,217.8159526443482,176.1856260204315,165.734731502533,193.0920477771759,191.0120236682892,160.5366449642182,165.4354525089265,172.0160270500184,111.9043052482605,182.9178998947144,198.4857194232941,165.8505724620819,164.0816167640686,183.431044216156,178.3703242588043,166.8244445610047,233.7849170303345,200.437901210785,166.2358214759827,200.8662677001954,211.3945352363587,205.8576155471802,176.3785883140564,218.7188100624085,217.278196258545,146.4260239028931,182.876769399643,198.1434554195404,180.6501061344147,172.3198363590241,186.027047405243,214.1142080879212,227.2311932182312,249.8055694007875,191.4335284519196,176.8794663429261,180.235184879303,216.8981531333924,190.3138517284394,145.8859228801728,230.5218712615967,207.9059804344178,161.1435615158081,198.9178479290009,204.0207845973969,199.5044840335846,175.8991564273834,169.0430141925812,194.6130018234253,144.1894919681549,224.3374343872071)

# Update existing rows 2..101 (data index 0..99) and add new rows 102..166 (data index 100..164)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    if ($row -gt 101) {
        $ws.Cells.Item($row, 1).Value = $i
    }
}

# Apply the same style as the existing index column (A) cells to the newly added ones
$ws.Range("A2").Copy()
$ws.Range("A102:A166").PasteSpecial(-4122)
